$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sub-header row (Hiver / Été / Année), shifting data rows up.
$ws.Rows(2).Delete()

# New column headers in row 1 (idx / idx2 / Name / Date Start / Date End use
# the plain default style; the unit headers reuse the existing 9pt font).
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Match the selection left by the author after editing the header row.
$ws.Range("A2:K2").Select()
